# Daily attendance processing - 2026-02-07 08:00:05 UTC
# Correct the "Recorded By" column (G) formatting: swap from
# "Miss Dina Nasr, Administrator" to "Administrator, Miss Dina Nasr"
# for every populated session row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "Miss Dina Nasr, Administrator"
$newText = "Administrator, Miss Dina Nasr"

$lastRow = $ws.UsedRange.Rows.Count
$replaced = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
        $replaced++
    }
}

Write-Output "Updated $replaced 'Recorded By' cell(s) in column G."
